# Insert a new weekly price record as row 19, pushing the existing
# rows 19:25 down to 20:26 (dimension grows from A1:R25 to A1:R26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 19, shifting rows 19:25 -> 20:26.
$ws.Rows.Item(19).EntireRow.Insert()

# Populate the newly-inserted row 19 with the new observation.
$ws.Cells.Item(19, 1).Value = 6
$ws.Cells.Item(19, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 44754
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = 100112035
$ws.Cells.Item(19, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 19000
$ws.Cells.Item(19, 13).Value = 18133
$ws.Cells.Item(19, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(19, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(19, 16).Value = 1209
$ws.Cells.Item(19, 17).Value = 15
$ws.Cells.Item(19, 18).Value = "Hortaliza"
